$d = $word.ActiveDocument

$replacements = @(
    @("33÷3=", "93÷6="),
    @("32÷4=", "85÷3="),
    @("17÷5=", "40÷7="),
    @("31÷9=", "71÷9="),
    @("90÷2=", "28÷5="),
    @("85÷4=", "19÷2="),
    @("31÷5=", "90÷5="),
    @("61÷2=", "88÷8="),
    @("26÷9=", "82÷7="),
    @("98÷9=", "51÷3="),
    @("99÷4=", "81÷4="),
    @("27÷5=", "45÷8="),
    @("97÷8=", "16÷2="),
    @("16÷7=", "58÷2="),
    @("46÷2=", "22÷8="),
    @("22÷3=", "75÷4="),
    @("84÷2=", "57÷4="),
    @("62÷3=", "74÷2="),
    @("35÷3=", "27÷5="),
    @("43÷3=", "18÷5="),
    @("70÷2=", "48÷2="),
    @("49÷3=", "15÷6="),
    @("56÷5=", "88÷8="),
    @("85÷7=", "45÷2="),
    @("30÷4=", "83÷5="),
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]
    $result = $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output "WARNING: replacement failed for $find -> $replace"
    }
}

Write-Output "Done."